$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.705.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.260.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.46%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.75'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.02'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.40%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.479'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.12'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.17'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0794'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.83%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.56'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.615.50'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.13'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.275.00'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.758'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.626.49'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.36'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +10.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0900'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.88'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.56'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.52'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.57'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.23%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.03'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.47'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.10%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.47'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +8.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '160.63'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.13'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0742'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.00'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.01%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.52'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.73%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.79'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.92'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.059.46'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.11'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.46%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.04'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.23%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.85'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '72.43'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.77%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.04%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.88%  '
